$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# yellow fill color used for the highlighted cells (BGR packed value -> RGB FFFF00)
$yellow = 65535

# --- Row 31: get the 0x48 ntc temperature ---
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "ntc_temp_get.sh"
$ws.Range("C31").Value = "./ntc_temp_get.sh <ntc_addr>"
$ws.Range("D31").Value = "adb shell /etc/factory-test/lvp15/ntc_temp_get.sh 0x48"
$ws.Range("E31").Value = "get the 0x48 ntc temperature "
$ws.Range("A31:E31").Interior.Color = $yellow

# --- Row 32: get the 0x49 ntc temperature ---
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "ntc_temp_get.sh"
$ws.Range("C32").Value = "./ntc_temp_get.sh <ntc_addr>"
$ws.Range("D32").Value = "adb shell /etc/factory-test/lvp15/ntc_temp_get.sh 0x49"
$ws.Range("E32").Value = "get the 0x49 ntc temperature"
$ws.Range("A32:E32").Interior.Color = $yellow

# --- Row 33: get the 0x4a ntc temperature ---
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "ntc_temp_get.sh"
$ws.Range("C33").Value = "./ntc_temp_get.sh <ntc_addr>"
$ws.Range("D33").Value = "adb shell /etc/factory-test/lvp15/ntc_temp_get.sh 0x4a"
$ws.Range("E33").Value = "get the 0x4a ntc temperature"
$ws.Range("A33:E33").Interior.Color = $yellow

# --- Row 34: get the 0x4b ntc temperature ---
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "ntc_temp_get.sh"
$ws.Range("C34").Value = "./ntc_temp_get.sh <ntc_addr>"
$ws.Range("D34").Value = "adb shell /etc/factory-test/lvp15/ntc_temp_get.sh 0x4b"
$ws.Range("E34").Value = "get the 0x4b ntc temperature"
$ws.Range("A34:E34").Interior.Color = $yellow

# --- Row 35: mute amp ---
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "set_amp_mute.sh"
$ws.Range("C35").Value = "./set_amp_mute.sh"
$ws.Range("D35").Value = "adb shell /etc/factory-test/lvp15/set_amp_mute.sh"
$ws.Range("E35").Value = "mute amp"
$ws.Range("A35").Interior.Color = $yellow
$ws.Range("E35").Interior.Color = $yellow

# --- Row 36: unmute amp ---
$ws.Range("A36").Value = 33
$ws.Range("B36").Value = "set_amp_unmute.sh"
$ws.Range("C36").Value = "./set_amp_unmute.sh"
$ws.Range("D36").Value = "adb shell /etc/factory-test/lvp15/set_amp_unmute.sh"
$ws.Range("E36").Value = "unmute amp"
$ws.Range("A36").Interior.Color = $yellow
$ws.Range("E36").Interior.Color = $yellow

# --- update the view: scroll position and active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C28").Select()
